$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for columns I (Celular_Deputado), J (Celular_Assessoria) and
# K (Nome_assessor) for data rows 2-22. Several rows previously held the
# placeholder "xxxxx" (or other stale data) which is now replaced with the
# real phone numbers / assessor names.
$data = @(
    @("071 99961-2530", "061 99981-9683", "Fred"),
    @("061 99654-9496", "061 99994-8910", "Hermes"),
    @("061 99102-5662", "061 99985-7345", "Noeli"),
    @("021 99386-4728", "061 99519-1440", "Marivaldo"),
    @("031 99984-1123", "061 99985-8537", "Rafael"),
    @("031 99862-1588", "061 98104-0205", "Débora"),
    @("082 99646-1515", "082 99971-0612", "Renato"),
    @("061 99815-0833", "021 98105-1277", "Maria"),
    @("061 99623-9970", "061 99103-0773", "Rodrigo"),
    @("051 99591-4462", "051 99463-7344", "Conrado"),
    @("061 98124-1234", "061 98215-2443", "Daniel"),
    @("073 98834-8409", "061 99612-9930", "Cristiano"),
    @("081 99960-6635", "061 98116-1635", "Fabiano"),
    @("098 99105-0044", "061 99943-0014", "Marcos"),
    @("011 97317-1777", "061 98570-0340", "Liliene"),
    @("061 99175-8139", "061 99131-0375", "Monalisa"),
    @("021 98187-4793", "022 98112-3198", "Emanuel"),
    @("061 98157-1312", "061 99820-3010", "Teresa"),
    @("061 98118-4600", "061 98118-1857", "Solange"),
    @("049 98831-2168", "061 99403-1487", "Flávia"),
    @("051 99774-1717", "061 99108-0813", "Marina")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 9).Value = $entry[0]
    $ws.Cells.Item($row, 10).Value = $entry[1]
    $ws.Cells.Item($row, 11).Value = $entry[2]
    $row++
}

# Update the active selection left by the editor when the file was saved.
$ws.Range("K25").Select()

$wb.Save()
